$newOrder = @(
    "A 61919-2021",
    "A 10000-2022",
    "A 39620-2023",
    "A 39613-2023",
    "A 51848-2024",
    "A 11714-2021",
    "A 22400-2025",
    "A 34068-2022",
    "A 64357-2023",
    "A 3677-2026",
    "A 33751-2024",
    "A 15452-2023",
    "A 51189-2022",
    "A 64345-2023",
    "A 39716-2023",
    "A 9750-2022",
    "A 53462-2021",
    "A 11477-2021",
    "A 74167-2021",
    "A 11459-2021",
    "A 35953-2021",
    "A 11713-2021",
    "A 53090-2021",
    "A 31697-2022",
    "A 11197-2021",
    "A 9767-2022",
    "A 65340-2021",
    "A 55112-2021",
    "A 25459-2023",
    "A 7214-2022",
    "A 11285-2023",
    "A 25452-2023",
    "A 49592-2023",
    "A 60865-2024",
    "A 40643-2024",
    "A 6826-2025",
    "A 16734-2025",
    "A 51770-2023",
    "A 41949-2025",
    "A 41963-2025",
    "A 41951-2025",
    "A 41957-2025",
    "A 31684-2022",
    "A 41841-2025",
    "A 35966-2024",
    "A 61589-2024",
    "A 46757-2023",
    "A 3890-2025",
    "A 43158-2024",
    "A 24946-2024",
    "A 21513-2025",
    "A 9856-2024",
    "A 37392-2024",
    "A 19031-2024",
    "A 35961-2024",
    "A 22402-2025",
    "A 48375-2024",
    "A 24636-2023",
    "A 10790-2023",
    "A 24105-2025",
    "A 24967-2024",
    "A 862-2022",
    "A 60934-2023",
    "A 41693-2022",
    "A 26256-2025",
    "A 12931-2022",
    "A 74170-2021",
    "A 29193-2025",
    "A 47392-2025",
    "A 28752-2024",
    "A 41702-2023",
    "A 36179-2022",
    "A 17286-2022",
    "A 28681-2025",
    "A 49916-2025",
    "A 42572-2024",
    "A 45233-2024",
    "A 59957-2023",
    "A 11718-2021",
    "A 12900-2022",
    "A 31695-2022",
    "A 17320-2023",
    "A 5142-2024",
    "A 20957-2023",
    "A 8673-2025",
    "A 59165-2022",
    "A 28851-2022",
    "A 9120-2025",
    "A 61937-2021",
    "A 57871-2025",
    "A 57879-2025",
    "A 57881-2025",
    "A 45319-2023",
    "A 19288-2025",
    "A 19294-2025",
    "A 19311-2025",
    "A 2663-2026",
    "A 2675-2026",
    "A 2628-2026",
    "A 60935-2023",
    "A 32420-2023",
    "A 32426-2023",
    "A 15428-2021",
    "A 3673-2026",
    "A 4788-2026",
    "A 5175-2026",
    "A 20375-2023",
    "A 38276-2023",
    "A 2496-2024",
    "A 60949-2025",
    "A 59168-2022",
    "A 60969-2025",
    "A 72530-2021",
    "A 6893-2026",
    "A 7229-2026",
    "A 7833-2026",
    "A 20960-2023",
    "A 32424-2023",
    "A 61576-2024",
    "A 58619-2023",
    "A 43700-2023",
    "A 56796-2022",
    "A 36181-2022",
    "A 60922-2023",
    "A 273-2025",
    "A 38465-2023"
)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 127
$numRows = $lastRow - $firstRow + 1
$numCols = 26  # A..Z

$srcRange = $ws.Range("A" + $firstRow + ":Z" + $lastRow)
$data = $srcRange.Formula()

# Map label (column A) -> old row index (1-based within $data)
$labelToIndex = @{}
for ($r = 1; $r -le $numRows; $r++) {
    $label = $data[$r,1]
    $labelToIndex[$label] = $r
}

$newData = New-Object 'object[,]' $numRows,$numCols

for ($i = 0; $i -lt $numRows; $i++) {
    $label = $newOrder[$i]
    $oldR = $labelToIndex[$label]
    for ($c = 1; $c -le $numCols; $c++) {
        $val = $data[$oldR,$c]
        if ($c -eq 3) {
            $val = 46064
        }
        $newData[$i, $c-1] = $val
    }
}

$dstRange = $ws.Range("A" + $firstRow + ":Z" + $lastRow)
$dstRange.Formula = $newData
